# Fill in the "Rodada 28/29/30" columns (AC:AE) for rows 2-21 with the
# values that were entered in the source workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("AC2").Value = 92.84
$ws.Range("AD2").Value = 108.53
$ws.Range("AE2").Value = 23.99
$ws.Range("AC3").Value = 116.54
$ws.Range("AD3").Value = 108.58
$ws.Range("AE3").Value = 35.39
$ws.Range("AC4").Value = 127.14
$ws.Range("AD4").Value = 112.53
$ws.Range("AE4").Value = 64.19
$ws.Range("AC5").Value = 133.64
$ws.Range("AD5").Value = 122.58
$ws.Range("AE5").Value = 55.79
$ws.Range("AC6").Value = 134.94
$ws.Range("AD6").Value = 74.79
$ws.Range("AE6").Value = 37.59
$ws.Range("AC7").Value = 90.96
$ws.Range("AD7").Value = 70.18
$ws.Range("AE7").Value = 56
$ws.Range("AC8").Value = 131.04
$ws.Range("AD8").Value = 109.43
$ws.Range("AE8").Value = 57.3
$ws.Range("AC9").Value = 133.14
$ws.Range("AD9").Value = 92.88
$ws.Range("AE9").Value = 54.79
$ws.Range("AC10").Value = 107.06
$ws.Range("AD10").Value = 88.43
$ws.Range("AE10").Value = 56.04
$ws.Range("AC11").Value = 132.86
$ws.Range("AD11").Value = 122.6
$ws.Range("AE11").Value = 27.69
$ws.Range("AC12").Value = 85.18
$ws.Range("AD12").Value = 93.04
$ws.Range("AE12").Value = 31.93
$ws.Range("AC13").Value = 92.55
$ws.Range("AD13").Value = 77.06
$ws.Range("AE13").Value = 28.69
$ws.Range("AC14").Value = 95.44
$ws.Range("AD14").Value = 41.59
$ws.Range("AE14").Value = 41.14
$ws.Range("AC15").Value = 109.44
$ws.Range("AD15").Value = 91.43
$ws.Range("AE15").Value = 37.44
$ws.Range("AC16").Value = 111.94
$ws.Range("AD16").Value = 115.69
$ws.Range("AE16").Value = 43.14
$ws.Range("AC17").Value = 54.65
$ws.Range("AD17").Value = 47.71
$ws.Range("AE17").Value = 49.69
$ws.Range("AC18").Value = 49.98
$ws.Range("AD18").Value = 60.15
$ws.Range("AE18").Value = 19.25
$ws.Range("AC19").Value = 87.38
$ws.Range("AD19").Value = 68.24
$ws.Range("AE19").Value = 50.61
$ws.Range("AC20").Value = 26.55
$ws.Range("AD20").Value = 25.01
$ws.Range("AE20").Value = 57.27
$ws.Range("AC21").Value = 31.44
$ws.Range("AD21").Value = 14.09
$ws.Range("AE21").Value = 40.16

# AE2 ends up with a slightly different cell style than the rest of the
# AE column (border only on left/right, explicit "no fill") - reproduce
# that by copying the border layout from a cell that already uses it
# (U13) and then explicitly clearing the fill on AE2.
$ws.Range("U13").Copy() | Out-Null
$ws.Range("AE2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws.Range("AE2").Interior.ColorIndex = -4142

# Column AN (40) used to be auto-sized ("best fit"); it is now a fixed,
# wider custom width.
$ws.Columns.Item(40).ColumnWidth = 10.6

# The sheet no longer has a frozen/scrolled "top-left" cell, and the
# active selection moved from AF35 to AF28.
$ws.Range("AF28").Select()
